# Updated Figure for Paths
# Update the B2:E2 values ("Time taken", "distance offset", "distance traveled",
# "area between paths") on each "Trial No.N" worksheet to the new figures.

$wb = $excel.ActiveWorkbook

$updates = @{
    "Trial No.1"  = @(50.812, 0.2593528856428975, 35.37482070922852, 1.826734662055969)
    "Trial No.2"  = @(50.801, 0.2081243245494995, 35.7010612487793, 1.96355926990509)
    "Trial No.3"  = @(51.018, 0.28305833119626, 35.61456298828125, 2.108329296112061)
    "Trial No.4"  = @(51.41, 0.1756432299417202, 35.79704666137695, 1.762745022773743)
    "Trial No.5"  = @(51.208, 0.3612406802292721, 35.70381164550781, 1.866026639938354)
    "Trial No.6"  = @(51.013, 0.3700479789610016, 35.57451629638672, 2.190979242324829)
    "Trial No.7"  = @(84.806, 0.2400655657575073, 36.01694488525391, 1.417576432228088)
    "Trial No.8"  = @(51.027, 0.3040378721061801, 35.58433151245117, 2.12936544418335)
    "Trial No.9"  = @(51.209, 0.2428398027929857, 35.55540084838867, 1.886032342910767)
    "Trial No.10" = @(51.002, 0.2317258327108997, 35.46334838867188, 2.043991804122925)
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $vals = $updates[$sheetName]

    $ws.Range("B2").Value = $vals[0]
    $ws.Range("C2").Value = $vals[1]
    $ws.Range("D2").Value = $vals[2]
    $ws.Range("E2").Value = $vals[3]
}
